$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every data row (2..131): the "M" column comment text actually belongs
# in the "S" column. Move the value (and its cell format) from M to S, and
# blank out M (keeping M's own formatting, which does not change).
for ($r = 2; $r -le 131; $r++) {
    $mCell = $ws.Range("M$r")
    $sCell = $ws.Range("S$r")

    $mValue = $mCell.Value2

    # Bring S's formatting in line with M's (this also normalizes the
    # handful of rows whose S cell already had a slightly different style
    # even though M was empty).
    $mCell.Copy()
    $sCell.PasteSpecial(-4122)

    if ($mValue -ne $null) {
        $sCell.Value2 = $mValue
    }

    $mCell.ClearContents()
}

$excel.CutCopyMode = 0

# Refresh the saved view/selection state as recorded by the author.
$ws.Range("N15").Select()
$ws.Application.ActiveWindow.ScrollColumn = 9
